$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(10002, 110021),
    @(10003, 110022),
    @(10004, 110023),
    @(10005, 110024),
    @(10006, 110025),
    @(10007, 110026),
    @(10008, 110027),
    @(10009, 110028),
    @(10010, 110029)
)

$row = 22
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $row = $row + 1
}

[void]$ws.Rows("31:1048576").Select()

$ws.PageSetup.Orientation = 1
